# Automatic BRVM update (GitHub Actions) - applies refreshed market data
# to the recommendations worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: sector/index summary rows (only "Jours en Baisse" (C) and
# --- "Variation Totale (%)" (D) columns were refreshed) ---

$ws.Range("C2").Value = 130
$ws.Range("D2").Value = 57094.79

$ws.Range("C3").Value = 52
$ws.Range("D3").Value = 47235

$ws.Range("C4").Value = 65
$ws.Range("D4").Value = 44238.71

$ws.Range("C5").Value = 53
$ws.Range("D5").Value = 44165

$ws.Range("C6").Value = 60
$ws.Range("D6").Value = 43350

$ws.Range("C8").Value = 65
$ws.Range("D8").Value = 38490

$ws.Range("C9").Value = 65
$ws.Range("D9").Value = 37480

$ws.Range("C10").Value = 65
$ws.Range("D10").Value = 32255

$ws.Range("C11").Value = 65
$ws.Range("D11").Value = 26050

$ws.Range("C12").Value = 65
$ws.Range("D12").Value = 23922.88

$ws.Range("C13").Value = 65
$ws.Range("D13").Value = 21474.45

$ws.Range("C14").Value = 65
$ws.Range("D14").Value = 13971.55

$ws.Range("C15").Value = 65
$ws.Range("D15").Value = 9556.809999999999

$ws.Range("C16").Value = 65
$ws.Range("D16").Value = 8555.690000000001

$ws.Range("C17").Value = 65
$ws.Range("D17").Value = 7868.62

$ws.Range("C18").Value = 65
$ws.Range("D18").Value = 7325.82

$ws.Range("C19").Value = 65
$ws.Range("D19").Value = 7254.25

$ws.Range("C20").Value = 65
$ws.Range("D20").Value = 7096.64

$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 6995.6

$ws.Range("C22").Value = 65
$ws.Range("D22").Value = 6875.07

$ws.Range("C23").Value = 65
$ws.Range("D23").Value = 6515.08

$ws.Range("C24").Value = 65
$ws.Range("D24").Value = 6455.58

# --- Section 2: stock rows 40-64 were re-sorted by "Variation Totale (%)"
# --- (descending) with refreshed B/C/D/E/F values. Row 54 ("TOTAL") stays
# --- in place as a separator between "Achat"/"Observer" and "Vente" blocks.

$data = @(
    @{ Row = 40; A = "AFRICA GLOBAL LOGISTICS CI (SDSC)";           B = 5;  C = 0;  D = 15.74;              E = 4.1;   F = "🟢 Achat" },
    @{ Row = 41; A = "UNIWAX CI (UNXC)";                            B = 12; C = 9;  D = 14.12;              E = 2.53;  F = "🟢 Achat" },
    @{ Row = 42; A = "BANK OF AFRICA NG (BOAN)";                    B = 13; C = 12; D = 13.55;              E = 4.49;  F = "🟢 Achat" },
    @{ Row = 43; A = "SUCRIVOIRE (SCRC)";                           B = 8;  C = 9;  D = 11.75;              E = -1.01; F = "🟢 Achat" },
    @{ Row = 44; A = "TRACTAFRIC MOTORS CI (PRSC)";                 B = 10; C = 10; D = 10.12;              E = -6.49; F = "🟢 Achat" },
    @{ Row = 45; A = "CFAO MOTORS CI (CFAC)";                       B = 7;  C = 9;  D = 9.92;               E = -6.56; F = "🟢 Achat" },
    @{ Row = 46; A = "SMB CI (SMBC)";                               B = 11; C = 12; D = 9.83;               E = -1.64; F = "🟢 Achat" },
    @{ Row = 47; A = "ONATEL BF (ONTBF)";                           B = 5;  C = 7;  D = 9.779999999999999;  E = -1.28; F = "🟢 Achat" },
    @{ Row = 48; A = "SODE CI (SDCC)";                              B = 5;  C = 7;  D = 5.85;               E = -1.83; F = "🟢 Achat" },
    @{ Row = 49; A = "VIVO ENERGY CI (SHEC)";                       B = 4;  C = 4;  D = 5.55;               E = 2.63;  F = "🟢 Achat" },
    @{ Row = 50; A = "ECOBANK COTE D''IVOIRE (ECOC)";               B = 5;  C = 4;  D = 4.77;               E = 3.94;  F = "🟡 Observer" },
    @{ Row = 51; A = "SETAO CI (STAC)";                             B = 14; C = 11; D = 4.1;                E = 1.83;  F = "🟡 Observer" },
    @{ Row = 52; A = "BICI CI (BICC)";                              B = 1;  C = 1;  D = 3.7;                E = -1.22; F = "🟡 Observer" },
    @{ Row = 53; A = "UNILEVER CI (UNLC)";                          B = 7;  C = 7;  D = 2.53;               E = -7.46; F = "🟡 Observer" },
    @{ Row = 55; A = "SONATEL SN (SNTS)";                           B = 2;  C = 3;  D = -3.43;              E = 0.8;   F = "🟡 Observer" },
    @{ Row = 56; A = "SOGB CI (SOGC)";                              B = 6;  C = 5;  D = -3.76;              E = 2.78;  F = "🟡 Observer" },
    @{ Row = 57; A = "TOTALENERGIES MARKETING SN (TTLS)";           B = 10; C = 14; D = -4.09;              E = 2.34;  F = "🟡 Observer" },
    @{ Row = 58; A = "NEI-CEDA CI (NEIC)";                          B = 6;  C = 8;  D = -5.59;              E = 5.88;  F = "🔴 Vente" },
    @{ Row = 59; A = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)";       B = 6;  C = 9;  D = -6.92;              E = -1.35; F = "🔴 Vente" },
    @{ Row = 60; A = "SOLIBRA CI (SLBC)";                           B = 11; C = 13; D = -8.32;              E = 3.91;  F = "🔴 Vente" },
    @{ Row = 61; A = "LOTERIE NATIONALE DU BENIN (LNBB)";           B = 0;  C = 5;  D = -9.130000000000001; E = -2.63; F = "🔴 Vente" },
    @{ Row = 62; A = "ORANGE COTE D'IVOIRE (ORAC)";                 B = 9;  C = 12; D = -11.95;             E = 2.76;  F = "🔴 Vente" },
    @{ Row = 63; A = "ORAGROUP TOGO (ORGT)";                        B = 6;  C = 8;  D = -12.37;             E = 0.31;  F = "🔴 Vente" },
    @{ Row = 64; A = "CORIS BANK INTERNATIONAL (CBIBF)";            B = 3;  C = 9;  D = -20.64;             E = -2.78; F = "🔴 Vente" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
}

# Row 54 ("TOTAL") keeps its position; only "Jours en Baisse" (C) changes.
$ws.Range("C54").Value = 64
